# Update the "cryptos" price/volume table (columns D and E) for the
# rows that changed in this run's snapshot. Price values that would
# otherwise be auto-parsed as numbers by Excel are written with a
# leading apostrophe so they stay text cells (matching the source
# data, which stores every Price/Volume cell as a string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.462.31'
$ws.Range("E2").Value = '  -3.00%  '
$ws.Range("D3").Value = '1.741.98'
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''321.87'
$ws.Range("D6").Value = '''0.9995'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.4252'
$ws.Range("E7").Value = '  -8.67%  '
$ws.Range("D8").Value = '''0.3594'
$ws.Range("E8").Value = '  -5.77%  '
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").Value = '''0.07413'
$ws.Range("E10").Value = '  -3.24%  '
$ws.Range("D11").Value = '''1.112'
$ws.Range("E11").Value = '  -4.17%  '
$ws.Range("D12").Value = '''0.9999'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").Value = '''21.39'
$ws.Range("E13").Value = '  -5.14%  '
$ws.Range("D14").Value = '''6.086'
$ws.Range("D15").Value = '''7.192'
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("D16").Value = '1.737.35'
$ws.Range("D18").Value = '''86.99'
$ws.Range("E18").Value = '  +6.11%  '
$ws.Range("D19").Value = '''0.06272'
$ws.Range("E19").Value = '  -6.62%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '''16.85'
$ws.Range("E21").Value = '  -3.86%  '
$ws.Range("D22").Value = '''6.102'
$ws.Range("E22").Value = '  -5.20%  '
$ws.Range("D23").Value = '''0.5227'
$ws.Range("E23").Value = '  -6.63%  '
$ws.Range("D24").Value = '27.482.53'
$ws.Range("E24").Value = '  -2.91%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").Value = '''2.318'
$ws.Range("E26").Value = '  -4.04%  '
$ws.Range("D27").Value = '''20.37'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '''150.92'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("D29").Value = '''2.345'
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("D30").Value = '1.936.06'
$ws.Range("E30").Value = '  -3.90%  '
$ws.Range("D31").Value = '''1.218'
$ws.Range("E31").Value = '  -3.52%  '
$ws.Range("D32").Value = '''126.57'
$ws.Range("E32").Value = '  -5.34%  '
$ws.Range("D33").Value = '''5.680'
$ws.Range("E33").Value = '  -3.41%  '
$ws.Range("D34").Value = '''0.09136'
$ws.Range("E34").Value = '  -4.75%  '
$ws.Range("E35").Value = '  -9.13%  '
$ws.Range("D36").Value = '''12.70'
$ws.Range("E36").Value = '  +4.52%  '
$ws.Range("D37").Value = '''0.02290'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").Value = '''0.2134'
$ws.Range("E38").Value = '  -6.47%  '
$ws.Range("D39").Value = '''5.068'
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("D40").Value = '''0.06068'
$ws.Range("E40").Value = '  -5.18%  '
$ws.Range("D41").Value = '''0.6397'
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("D42").Value = '''1.192'
$ws.Range("E42").Value = '  -3.96%  '
$ws.Range("D43").Value = '''1.416'
$ws.Range("E43").Value = '  -5.22%  '
$ws.Range("D44").Value = '''0.9993'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '''7.895'
$ws.Range("E45").Value = '  -5.30%  '
$ws.Range("D46").Value = '''13.75'
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("D47").Value = '''3.716'
$ws.Range("E47").Value = '  -3.77%  '
$ws.Range("D48").Value = '''0.5875'
$ws.Range("E48").Value = '  -4.69%  '
$ws.Range("D49").Value = '''125.33'
$ws.Range("E49").Value = '  -4.43%  '
$ws.Range("D50").Value = '''1.952'
$ws.Range("D51").Value = '''0.06854'
$ws.Range("E51").Value = '  -4.28%  '
